$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.213.82"
$ws.Range("E2").Value = "  +1.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.496.47"
$ws.Range("E3").Value = "  +1.55%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.35"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.10"
$ws.Range("E6").Value = "  +3.22%  "

$ws.Range("E7").Value = "  +1.18%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.13"
$ws.Range("E10").Value = "  +8.12%  "

$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.42"
$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.20"
$ws.Range("E14").Value = "  +1.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.887.12"
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.499.20"
$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("E17").Value = "  +0.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.161.77"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("E20").Value = "  +3.06%  "

$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("E22").Value = "  +13.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.71"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.70"
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("E25").Value = "  +2.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.03"
$ws.Range("E26").Value = "  -0.58%  "

$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  +0.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.09"
$ws.Range("E29").Value = "  +3.08%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("E30").Value = "  +8.97%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.34"
$ws.Range("E31").Value = "  +1.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.08"
$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.99"
$ws.Range("E33").Value = "  +0.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.42"
$ws.Range("E34").Value = "  +1.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0788"
$ws.Range("E35").Value = "  +2.93%  "

$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("E37").Value = "  +3.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.70"
$ws.Range("E38").Value = "  +2.07%  "

$ws.Range("E39").Value = "  +0.79%  "

$ws.Range("E40").Value = "  +0.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.11"
$ws.Range("E41").Value = "  -0.56%  "

$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.17"
$ws.Range("E43").Value = "  +0.72%  "

$ws.Range("E44").Value = "  +2.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.993.86"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.05"
$ws.Range("E46").Value = "  +1.84%  "

$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("E48").Value = "  -3.90%  "

$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("E50").Value = "  +0.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.78"
$ws.Range("E51").Value = "  +3.74%  "
